# Update "Сводная таблица" (summary table): rename the algorithm label
# from "ММП v2" to "ММП" and refresh the "t" / "t_max" timing columns
# that changed after fixing bugs / adding the dynamic algorithm.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($table, $rowIdx, $colIdx, $value) {
    # $rowIdx / $colIdx are 1-based Word COM indices.
    $cell = $table.Rows.Item($rowIdx).Cells.Item($colIdx)
    $cell.Range.Text = $value
}

# Table columns (1-based): 1=Алгоритм ... 12=t 13=t_max

Set-Cell $t 2 1  "ММП"
Set-Cell $t 2 12 "0.0038"
Set-Cell $t 2 13 "0.0575"

Set-Cell $t 3 13 "0.0022"

Set-Cell $t 4 1  "ММП"
Set-Cell $t 4 12 "0.0028"
Set-Cell $t 4 13 "0.0139"

Set-Cell $t 5 13 "0.0013"

Set-Cell $t 6 1  "ММП"
Set-Cell $t 6 12 "0.0031"
Set-Cell $t 6 13 "0.0152"

Set-Cell $t 7 12 "0.0003"
Set-Cell $t 7 13 "0.0028"

Set-Cell $t 8 1  "ММП"
Set-Cell $t 8 12 "0.0038"
Set-Cell $t 8 13 "0.0664"

Set-Cell $t 9 13 "0.0050"

Set-Cell $t 10 1  "ММП"
Set-Cell $t 10 12 "0.0036"
Set-Cell $t 10 13 "0.0427"

Set-Cell $t 11 13 "0.0010"

Set-Cell $t 12 1  "ММП"
Set-Cell $t 12 12 "0.0040"
Set-Cell $t 12 13 "0.0603"

Set-Cell $t 13 13 "0.0082"

Write-Host "Done updating table cells."
